$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated AgTests / AgPosit values (commit: Updated: st 08. 04. 2021)
$ws.Range("F265").Value = 16037
$ws.Range("G265").Value = 852
$ws.Range("F267").Value = 12806
$ws.Range("G267").Value = 764
$ws.Range("F268").Value = 13117
$ws.Range("F271").Value = 45776
$ws.Range("G271").Value = 1728
$ws.Range("F272").Value = 30775
$ws.Range("G272").Value = 1659
$ws.Range("F273").Value = 31688
$ws.Range("G273").Value = 1662
$ws.Range("F274").Value = 28101
$ws.Range("G274").Value = 1276
$ws.Range("F275").Value = 30347
$ws.Range("G275").Value = 1276
$ws.Range("F276").Value = 11344
$ws.Range("G276").Value = 379
$ws.Range("F278").Value = 30435
$ws.Range("G278").Value = 2097
$ws.Range("F279").Value = 42704
$ws.Range("G279").Value = 3031
$ws.Range("F280").Value = 34844
$ws.Range("G280").Value = 2320
$ws.Range("F281").Value = 46100
$ws.Range("G281").Value = 3168
$ws.Range("F282").Value = 46300
$ws.Range("G282").Value = 2754
$ws.Range("F286").Value = 54859
$ws.Range("G286").Value = 4278
$ws.Range("F287").Value = 58604
$ws.Range("G287").Value = 3716
$ws.Range("F288").Value = 59202
$ws.Range("G288").Value = 3971
$ws.Range("F289").Value = 63651
$ws.Range("G289").Value = 3681
$ws.Range("F292").Value = 82350
$ws.Range("G292").Value = 7262
$ws.Range("F293").Value = 82411
$ws.Range("G293").Value = 5760
$ws.Range("F294").Value = 93902
$ws.Range("G294").Value = 4942
$ws.Range("F295").Value = 17190
$ws.Range("G295").Value = 1030
$ws.Range("F297").Value = 2311
$ws.Range("F298").Value = 3210
$ws.Range("F299").Value = 65614
$ws.Range("G299").Value = 6874
$ws.Range("F300").Value = 72268
$ws.Range("G300").Value = 6962
$ws.Range("F301").Value = 72164
$ws.Range("G301").Value = 5682
$ws.Range("F362").Value = 228263
$ws.Range("G362").Value = 3176
$ws.Range("F363").Value = 186771
$ws.Range("G363").Value = 2749
$ws.Range("F364").Value = 167417
$ws.Range("G364").Value = 2463
$ws.Range("F365").Value = 183586
$ws.Range("G365").Value = 2383
$ws.Range("F366").Value = 338864
$ws.Range("G366").Value = 2837
$ws.Range("F367").Value = 765564
$ws.Range("G367").Value = 3917
$ws.Range("F368").Value = 345589
$ws.Range("G368").Value = 2294
$ws.Range("F369").Value = 233224
$ws.Range("G369").Value = 2590
$ws.Range("F370").Value = 181988
$ws.Range("G370").Value = 2035
$ws.Range("F371").Value = 159213
$ws.Range("G371").Value = 1948
$ws.Range("F372").Value = 178266
$ws.Range("G372").Value = 1848
$ws.Range("F373").Value = 348036
$ws.Range("G373").Value = 2370
$ws.Range("F374").Value = 770910
$ws.Range("G374").Value = 3418
$ws.Range("F376").Value = 220803
$ws.Range("G376").Value = 2223
$ws.Range("F377").Value = 176132
$ws.Range("G377").Value = 1819
$ws.Range("F378").Value = 156956
$ws.Range("G378").Value = 1544
$ws.Range("F379").Value = 178771
$ws.Range("G379").Value = 1605
$ws.Range("F380").Value = 343570
$ws.Range("G380").Value = 2008
$ws.Range("F381").Value = 743030
$ws.Range("G381").Value = 2682
$ws.Range("F383").Value = 220394
$ws.Range("G383").Value = 1756
$ws.Range("F384").Value = 171181
$ws.Range("G384").Value = 1506
$ws.Range("F385").Value = 150500
$ws.Range("G385").Value = 1400
$ws.Range("F386").Value = 182059
$ws.Range("G386").Value = 1355
$ws.Range("F387").Value = 350927
$ws.Range("G387").Value = 1666
$ws.Range("F388").Value = 718867
$ws.Range("G388").Value = 2159
$ws.Range("F389").Value = 350759
$ws.Range("G389").Value = 1298
$ws.Range("F390").Value = 218989
$ws.Range("G390").Value = 1504
$ws.Range("F391").Value = 175878
$ws.Range("G391").Value = 1202
$ws.Range("F392").Value = 216147
$ws.Range("G392").Value = 1193
$ws.Range("F393").Value = 293130
$ws.Range("G393").Value = 1173
$ws.Range("F394").Value = 160337
$ws.Range("G394").Value = 612
$ws.Range("F395").Value = 724832
$ws.Range("G395").Value = 1893
$ws.Range("F396").Value = 161802
$ws.Range("G396").Value = 541
$ws.Range("F397").Value = 103820
$ws.Range("G397").Value = 616
$ws.Range("F398").Value = 283373
$ws.Range("G398").Value = 1414

# Retracted AgTests/AgPosit figures for rows 334-361 (cells removed entirely)
$ws.Range("F334:G361").ClearContents()

